# Scheduled-runner market data refresh: updates the independently-sourced
# price columns (currentAveragePrice / NQ / HQ) on affected Leve rows across
# the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, then writes through the
# dependent LevePriceNQ/HQ + LeveProfitNQ/HQ columns (Price = unitPrice *
# LeveAmount; Profit = LeveGil - PriceNQ resp. -(PriceHQ + 2*LeveGil)).
# A Profit cell is only ever present when its corresponding Price is
# non-zero, so some cells are newly created while others are cleared here,
# matching upstream's "omit non-negative / undefined profit" convention.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3791.7646
$ws.Range("J17").Value = 4084.516
$ws.Range("L17").Value = 12253.548
$ws.Range("N17").Value = -12589.548

$ws.Range("H43").Value = 8749.1
$ws.Range("I43").Value = 5475.5
$ws.Range("J43").Value = 10931.5
$ws.Range("K43").Value = 5475.5
$ws.Range("L43").Value = 10931.5
$ws.Range("M43").Value = -5406.5
$ws.Range("N43").Value = -11069.5

$ws.Range("H80").Value = 2653.5
$ws.Range("J80").Value = 3927.2727
$ws.Range("L80").Value = 11781.8181
$ws.Range("N80").Value = -13777.8181

$ws.Range("H83").Value = 2653.5
$ws.Range("J83").Value = 3927.2727
$ws.Range("L83").Value = 35345.4543
$ws.Range("N83").Value = -45329.4543

$ws.Range("H86").Value = 5486
$ws.Range("J86").Value = 6999.75
$ws.Range("L86").Value = 6999.75
$ws.Range("N86").Value = -9245.75

$ws.Range("H89").Value = 5486
$ws.Range("J89").Value = 6999.75
$ws.Range("L89").Value = 34998.75
$ws.Range("N89").Value = -46230.75

$ws.Range("H113").Value = 7871.8184
$ws.Range("J113").Value = 8309
$ws.Range("L113").Value = 8309
$ws.Range("N113").Value = -14817

$ws.Range("H118").Value = 1332.3
$ws.Range("J118").Value = 1728
$ws.Range("L118").Value = 5184
$ws.Range("N118").Value = -8498

$ws.Range("H128").Value = 69705
$ws.Range("J128").Value = 69705
$ws.Range("L128").Value = 69705
$ws.Range("N128").Value = -79665

$ws.Range("H138").Value = 3396.125
$ws.Range("I138").Value = 3712.8572
$ws.Range("K138").Value = 11138.5716
$ws.Range("M138").Value = -5998.571599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 33916.363
$ws.Range("I43").Value = 34171
$ws.Range("J43").Value = 33859.777
$ws.Range("K43").Value = 34171
$ws.Range("L43").Value = 33859.777
$ws.Range("M43").Value = -33858
$ws.Range("N43").Value = -34485.777

$ws.Range("H112").Value = 54627.4
$ws.Range("J112").Value = 54627.4
$ws.Range("L112").Value = 54627.4
$ws.Range("N112").Value = -57581.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4225.4
$ws.Range("I86").Value = 4643.5713
$ws.Range("J86").Value = 3249.6667
$ws.Range("K86").Value = 4643.5713
$ws.Range("L86").Value = 3249.6667
$ws.Range("M86").Value = -3520.5713
$ws.Range("N86").Value = -5495.6667

$ws.Range("H89").Value = 4225.4
$ws.Range("I89").Value = 4643.5713
$ws.Range("J89").Value = 3249.6667
$ws.Range("K89").Value = 23217.8565
$ws.Range("L89").Value = 16248.3335
$ws.Range("M89").Value = -17601.8565
$ws.Range("N89").Value = -27480.3335

$ws.Range("H98").Value = 65000
$ws.Range("J98").Value = 65000
$ws.Range("L98").Value = 65000
$ws.Range("N98").Value = -70990

$ws.Range("H134").Value = 2290.0386
$ws.Range("I134").Value = 1271.375
$ws.Range("J134").Value = 14514
$ws.Range("K134").Value = 3814.125
$ws.Range("L134").Value = 43542
$ws.Range("M134").Value = -1279.125
$ws.Range("N134").Value = -48612

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 3000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 3000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 3000
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -3226

$ws.Range("H4").Value = 718571.4399999999
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 718571.4399999999
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 718571.4399999999
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -718795.4399999999

$ws.Range("H5").Value = 2743.5
$ws.Range("I5").Value = 56.5
$ws.Range("J5").Value = 4087
$ws.Range("K5").Value = 56.5
$ws.Range("L5").Value = 4087
$ws.Range("M5").Value = 55.5
$ws.Range("N5").Value = -4311

$ws.Range("H10").Value = 1544.7142
$ws.Range("I10").Value = 507.5
$ws.Range("J10").Value = 1959.6
$ws.Range("K10").Value = 507.5
$ws.Range("L10").Value = 1959.6
$ws.Range("M10").Value = -368.5
$ws.Range("N10").Value = -2237.6

$ws.Range("H11").Value = 493.5
$ws.Range("I11").Value = 687.5
$ws.Range("J11").Value = 396.5
$ws.Range("K11").Value = 687.5
$ws.Range("L11").Value = 396.5
$ws.Range("M11").Value = -547.5
$ws.Range("N11").Value = -676.5

$ws.Range("H12").Value = 1399.6666
$ws.Range("I12").Value = 1099.5
$ws.Range("J12").Value = 2000
$ws.Range("K12").Value = 1099.5
$ws.Range("L12").Value = 2000
$ws.Range("M12").Value = -929.5
$ws.Range("N12").Value = -2340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 17008
$ws.Range("J88").Value = 19016
$ws.Range("L88").Value = 57048
$ws.Range("N88").Value = -57904

$ws.Range("H91").Value = 17008
$ws.Range("J91").Value = 19016
$ws.Range("L91").Value = 57048
$ws.Range("N91").Value = -60012

$ws.Range("H93").Value = 19027
$ws.Range("J93").Value = 19027
$ws.Range("L93").Value = 57081
$ws.Range("N93").Value = -60825

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 2527.6667
$ws.Range("I43").Value = 343.75
$ws.Range("J43").Value = 19999
$ws.Range("K43").Value = 343.75
$ws.Range("L43").Value = 19999
$ws.Range("M43").Value = -192.75
$ws.Range("N43").Value = -20301

$ws.Range("H99").Value = 34282.75
$ws.Range("I99").Value = 36316.332
$ws.Range("J99").Value = 32249.166
$ws.Range("K99").Value = 36316.332
$ws.Range("L99").Value = 32249.166
$ws.Range("M99").Value = -34070.332
$ws.Range("N99").Value = -36741.166

$ws.Range("H126").Value = 5658.091
$ws.Range("I126").Value = 2690.2222
$ws.Range("J126").Value = 19013.5
$ws.Range("K126").Value = 8070.6666
$ws.Range("L126").Value = 57040.5
$ws.Range("M126").Value = -5600.6666
$ws.Range("N126").Value = -61980.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4139.5835
$ws.Range("I61").Value = 1067
$ws.Range("K61").Value = 1067
$ws.Range("M61").Value = -865

$ws.Range("H113").Value = 4139.5835
$ws.Range("I113").Value = 1067
$ws.Range("K113").Value = 1067
$ws.Range("M113").Value = 1103

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 34999.5
$ws.Range("J105").Value = 34999.5
$ws.Range("L105").Value = 34999.5
$ws.Range("N105").Value = -41987.5
